$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# A dash was inserted into the product name ("344MS-..." -> "344-MS-...").
# The value is duplicated on both the input and output sheets (cell B1 on
# each), so update it in both places.
$ws1.Range("B1").Value = "344-MS-EPP-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"
$ws2.Range("B1").Value = "344-MS-EPP-DB-SAR-REC-NON-RNI-CTRFD-DL-MD-TR-1-ONTIME"

# Reset the selection on the input sheet back to B1.
$null = $ws1.Range("B1").Select()

# Switch the active/selected tab to the output sheet, and select B1 there.
$null = $ws2.Activate()
$null = $ws2.Range("B1").Select()
